$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.898.76"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").Value = "1.634.30"
$ws.Range("E3").Value = "  +0.94%  "
$ws.Range("E4").Value = "  +0.87%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.27"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.520"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.89%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "28.76"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.54%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.260"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.51%  "
$ws.Range("E11").Value = "  -1.07%  "
$ws.Range("D12").Value = "1.868.56"
$ws.Range("E12").Value = "  +0.92%  "
$ws.Range("D13").Value = "1.638.97"
$ws.Range("E13").Value = "  +1.10%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.585"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +3.47%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "9.50"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +6.51%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.88"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.18%  "
$ws.Range("D17").Value = "29.918.41"
$ws.Range("E17").Value = "  +0.30%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "64.94"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.90%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "240.62"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.25%  "
$ws.Range("E20").Value = "  -0.33%  "
$ws.Range("E21").Value = "  +0.75%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.86"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +2.38%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.14"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.09%  "
$ws.Range("E24").Value = "  +3.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.72"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.80%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.52"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.67%  "
$ws.Range("E27").Value = "  -0.98%  "
$ws.Range("E28").Value = "  +0.59%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.77%  "
$ws.Range("E30").Value = "  +0.53%  "
$ws.Range("E31").Value = "  +1.41%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.38"
$ws.Range("D32").ClearFormats()
$ws.Range("E33").Value = "  -0.55%  "
$ws.Range("D34").Value = "1.423.62"
$ws.Range("E34").Value = "  +0.29%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.69"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +3.67%  "
$ws.Range("E36").Value = "  -0.85%  "
$ws.Range("E37").Value = "  -3.26%  "
$ws.Range("E38").Value = "  +1.32%  "
$ws.Range("E39").Value = "  +0.42%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "75.88"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +9.38%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.561"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.60%  "
$ws.Range("E42").Value = "  -0.27%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.833"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.81%  "
$ws.Range("E44").Value = "  +0.68%  "
$ws.Range("E45").Value = "  +0.92%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.01"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.06%  "
$ws.Range("D47").Value = "1.776.42"
$ws.Range("E47").Value = "  +0.89%  "
$ws.Range("E48").Value = "  -1.90%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "48.80"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -9.35%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "92.84"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +5.22%  "
$ws.Range("D51").Value = "0.0₆0110"
$ws.Range("E51").Value = "  +8.98%  "
